$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "treat_NA_as_group" row entirely; everything below shifts up by one.
$ws.Rows(17).Delete()

# project_name
$ws.Range("B3").Value = "iahp_panel_3"

# excluded_channels value (drop trailing ", NA")
$ws.Range("B5").Value = "B2M, DNA, Bead, LD, Live, Dead, ID, Cell-ID, Cell_ID"

# do_normalization: turn on
$ws.Range("B7").Value = 1

# norm_mode comment rewording
$ws.Range("C8").Value = "percentile or harmony (harmony looks weird, will fix)"

# anchor_ids
$ws.Range("B9").Value = "HC-050"

# do_database_injection: turn on
$ws.Range("B11").Value = 1

# do_analysis: turn off
$ws.Range("B13").Value = 0

# grouping_columns
$ws.Range("B15").Value = "visit, hc_pre, hc_post, pre_post"

# grouping_orders
$ws.Range("B16").Value = "HC, S1.1, S1.2; HC, S1.1; HC, S1.2; S1.1, S1.2"
$ws.Rows(16).RowHeight = 57.6

# data_subsets (row previously held "treat_NA_as_group" data, now repurposed)
$ws.Range("B17").Value = "Monos_and_DCs, CD4_T, CD8_T, B, TCRgd_T, NK"

# clustering_k
$ws.Range("B27").Value = 30

# ccp_delta_cutoff
$ws.Range("B29").Value = 0.007

# umap_n
$ws.Range("B31").Value = 15

# umap_min_dist
$ws.Range("B32").Value = 0.1
